$d = $word.ActiveDocument

# Fix the misspelling "Delivarable" -> "Deliverable"
$d.Content.Find.Execute("Delivarable", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Deliverable", 2)

# Move the _GoBack bookmark so it sits right after "Deliverable 1:" and
# before the inline picture (instead of after the picture).
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$findRange = $d.Content.Duplicate
$findRange.Find.Execute("Deliverable 1:", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$bmRange = $d.Range($findRange.End, $findRange.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
